# Updates cryptos list values (coin prices / 1h volume %) to match the latest
# scrape, per the commit "Updated cryptos list ... with GitHub Actions".
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). Row 1 is the header.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.421.43'
$ws.Range('E2').Value = '  +1.09%  '

$ws.Range('D3').Value = '1.851.25'
$ws.Range('E3').Value = '  +1.12%  '

$ws.Range('E4').Value = '  +0.11%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.07'
$ws.Range('E5').Value = '  +1.60%  '

$ws.Range('E6').Value = '  +0.13%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4744'
$ws.Range('E7').Value = '  +2.47%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2744'
$ws.Range('E8').Value = '  +1.90%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06322'
$ws.Range('E9').Value = '  +1.77%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '17.57'
$ws.Range('E10').Value = '  +9.85%  '

$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').Value = '1.851.88'
$ws.Range('E11').Value = '  +1.37%  '

$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07459'
$ws.Range('E12').Value = '  +1.60%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.947'
$ws.Range('E13').Value = '  +1.25%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '84.46'
$ws.Range('E14').Value = '  +1.95%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6220'
$ws.Range('E15').Value = '  +0.78%  '

$ws.Range('D16').Value = '30.382.07'

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '243.96'
$ws.Range('E17').Value = '  +7.90%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.001'

$ws.Range('E19').Value = '  +3.21%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007313'
$ws.Range('E20').Value = '  +1.05%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('E21').Value = '  +0.15%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.903'
$ws.Range('E22').Value = '  +1.97%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.899'
$ws.Range('E23').Value = '  +1.54%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '164.98'
$ws.Range('E24').Value = '  -0.33%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.082'
$ws.Range('E25').Value = '  -0.26%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '17.94'
$ws.Range('E26').Value = '  +1.29%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.864'

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1030'
$ws.Range('E28').Value = '  +1.41%  '

$ws.Range('E29').Value = '  -1.30%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.028'
$ws.Range('E30').Value = '  -0.16%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.813'
$ws.Range('E31').Value = '  +1.85%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.04828'
$ws.Range('E32').Value = '  +0.85%  '

$ws.Range('E33').Value = '  +0.03%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.6947'
$ws.Range('E34').Value = '  -0.08%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.700'
$ws.Range('E35').Value = '  +0.78%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.01897'
$ws.Range('E36').Value = '  +4.71%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.679'
$ws.Range('E37').Value = '  +2.90%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.998'
$ws.Range('E38').Value = '  +4.85%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.8728'
$ws.Range('E39').Value = '  -1.80%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '106.26'
$ws.Range('E40').Value = '  +3.08%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.001'
$ws.Range('E41').Value = '  +0.16%  '

$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.504'
$ws.Range('E42').Value = '  +0.96%  '

$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4047'
$ws.Range('E43').Value = '  +1.68%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.137'
$ws.Range('E44').Value = '  +3.94%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '62.91'
$ws.Range('E45').Value = '  +6.75%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1195'
$ws.Range('E46').Value = '  +1.01%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '33.66'
$ws.Range('E47').Value = '  +3.82%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.544'
$ws.Range('E48').Value = '  +1.09%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.05513'
$ws.Range('E49').Value = '  -0.13%  '

$ws.Range('E50').Value = '  -0.13%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3668'
$ws.Range('E51').Value = '  +1.43%  '
